$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle row 63/64 (B cells become blank-but-styled, C/D switch from style 2 to style 9) ---
$ws.Range("A63").Copy()
$ws.Range("B63").PasteSpecial(-4122)
$ws.Range("C63").PasteSpecial(-4122)
$ws.Range("D63").PasteSpecial(-4122)

$ws.Range("A64").Copy()
$ws.Range("B64").PasteSpecial(-4122)
$ws.Range("C64").PasteSpecial(-4122)
$ws.Range("D64").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- New row 65: cwl_log_effect_loaded ---
$ws.Range("A65").Value = "cwl_log_effect_loaded"
$ws.Range("C65").Value = "loaded EffectSetting/{0}: {1} > {2}"
$ws.Range("D65").Value = "loaded EffectSetting/{0}: {1} > {2}"

$ws.Range("A63").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("B65").PasteSpecial(-4122)
$ws.Range("C65").PasteSpecial(-4122)
$ws.Range("D65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# D65 gets its own distinct highlighted font (same look, different family id) like the new font added upstream
$ws.Range("D65").Font.Family = 2

$ws.Range("A65:D65").RowHeight = 23.25

Write-Host "edit complete"
